$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (rows 2-7) so the shared string table can be rebuilt
# in the exact insertion order required by the target workbook.
[void]$ws.Range("A2:T7").ClearContents()

# --- Text columns, filled column-by-column (A, then B, then C, then D) ---
# so that new shared strings are interned in the required order:
# FAPs, ECs, sCs, Tslp, Crlf2
# Column A
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(10, 1).Value = "sCs"

# Column B
$ws.Cells.Item(2, 2).Value = "Tslp"
$ws.Cells.Item(3, 2).Value = "Tslp"
$ws.Cells.Item(4, 2).Value = "Tslp"
$ws.Cells.Item(5, 2).Value = "Tslp"
$ws.Cells.Item(6, 2).Value = "Tslp"
$ws.Cells.Item(7, 2).Value = "Tslp"
$ws.Cells.Item(8, 2).Value = "Tslp"
$ws.Cells.Item(9, 2).Value = "Tslp"
$ws.Cells.Item(10, 2).Value = "Tslp"

# Column C
$ws.Cells.Item(2, 3).Value = "Crlf2"
$ws.Cells.Item(3, 3).Value = "Crlf2"
$ws.Cells.Item(4, 3).Value = "Crlf2"
$ws.Cells.Item(5, 3).Value = "Crlf2"
$ws.Cells.Item(6, 3).Value = "Crlf2"
$ws.Cells.Item(7, 3).Value = "Crlf2"
$ws.Cells.Item(8, 3).Value = "Crlf2"
$ws.Cells.Item(9, 3).Value = "Crlf2"
$ws.Cells.Item(10, 3).Value = "Crlf2"

# Column D
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(3, 4).Value = "ECs"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(9, 4).Value = "ECs"
$ws.Cells.Item(10, 4).Value = "sCs"

# --- Numeric columns E through T ---
# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.578851
$ws.Cells.Item(2, 8).Value = 1.736553
$ws.Cells.Item(2, 9).Value = 0.1369952482678105
$ws.Cells.Item(2, 10).Value = 0.1369952482678105
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 10.422708
$ws.Cells.Item(2, 14).Value = 31.268124
$ws.Cells.Item(2, 15).Value = 0.8537134449068419
$ws.Cells.Item(2, 16).Value = 0.8537134449068419
$ws.Cells.Item(2, 17).Value = 6.033194948508
$ws.Cells.Item(2, 18).Value = 54.298754536572
$ws.Cells.Item(2, 19).Value = 0.1169546853345805
$ws.Cells.Item(2, 20).Value = 0.1169546853345805

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.578851
$ws.Cells.Item(3, 8).Value = 1.736553
$ws.Cells.Item(3, 9).Value = 0.1369952482678105
$ws.Cells.Item(3, 10).Value = 0.1369952482678105
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.048156333333333
$ws.Cells.Item(3, 14).Value = 3.144469
$ws.Cells.Item(3, 15).Value = 0.08585342255879412
$ws.Cells.Item(3, 16).Value = 0.08585342255879413
$ws.Cells.Item(3, 17).Value = 0.6067263417063333
$ws.Cells.Item(3, 18).Value = 5.460537075357
$ws.Cells.Item(3, 19).Value = 0.01176151093808324
$ws.Cells.Item(3, 20).Value = 0.01176151093808324

# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.578851
$ws.Cells.Item(4, 8).Value = 1.736553
$ws.Cells.Item(4, 9).Value = 0.1369952482678105
$ws.Cells.Item(4, 10).Value = 0.1369952482678105
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.7378083333333333
$ws.Cells.Item(4, 14).Value = 2.213425
$ws.Cells.Item(4, 15).Value = 0.06043313253436396
$ws.Cells.Item(4, 16).Value = 0.06043313253436396
$ws.Cells.Item(4, 17).Value = 0.4270810915583333
$ws.Cells.Item(4, 18).Value = 3.843729824025
$ws.Cells.Item(4, 19).Value = 0.008279051995146685
$ws.Cells.Item(4, 20).Value = 0.008279051995146685

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.879048666666666
$ws.Cells.Item(5, 8).Value = 8.637146
$ws.Cells.Item(5, 9).Value = 0.681377395677141
$ws.Cells.Item(5, 10).Value = 0.681377395677141
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 10.422708
$ws.Cells.Item(5, 14).Value = 31.268124
$ws.Cells.Item(5, 15).Value = 0.8537134449068419
$ws.Cells.Item(5, 16).Value = 0.8537134449068419
$ws.Cells.Item(5, 17).Value = 30.007483570456
$ws.Cells.Item(5, 18).Value = 270.067352134104
$ws.Cells.Item(5, 19).Value = 0.5817010437451843
$ws.Cells.Item(5, 20).Value = 0.5817010437451843

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.879048666666666
$ws.Cells.Item(6, 8).Value = 8.637146
$ws.Cells.Item(6, 9).Value = 0.681377395677141
$ws.Cells.Item(6, 10).Value = 0.681377395677141
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.048156333333333
$ws.Cells.Item(6, 14).Value = 3.144469
$ws.Cells.Item(6, 15).Value = 0.08585342255879412
$ws.Cells.Item(6, 16).Value = 0.08585342255879413
$ws.Cells.Item(6, 17).Value = 3.017693093941555
$ws.Cells.Item(6, 18).Value = 27.159237845474
$ws.Cells.Item(6, 19).Value = 0.05849858147308024
$ws.Cells.Item(6, 20).Value = 0.05849858147308025

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.879048666666666
$ws.Cells.Item(7, 8).Value = 8.637146
$ws.Cells.Item(7, 9).Value = 0.681377395677141
$ws.Cells.Item(7, 10).Value = 0.681377395677141
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.7378083333333333
$ws.Cells.Item(7, 14).Value = 2.213425
$ws.Cells.Item(7, 15).Value = 0.06043313253436396
$ws.Cells.Item(7, 16).Value = 0.06043313253436396
$ws.Cells.Item(7, 17).Value = 2.124186098338888
$ws.Cells.Item(7, 18).Value = 19.11767488505
$ws.Cells.Item(7, 19).Value = 0.04117777045887641
$ws.Cells.Item(7, 20).Value = 0.04117777045887641

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.7674366666666667
$ws.Cells.Item(8, 8).Value = 2.30231
$ws.Cells.Item(8, 9).Value = 0.1816273560550486
$ws.Cells.Item(8, 10).Value = 0.1816273560550485
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 10.422708
$ws.Cells.Item(8, 14).Value = 31.268124
$ws.Cells.Item(8, 15).Value = 0.8537134449068419
$ws.Cells.Item(8, 16).Value = 0.8537134449068419
$ws.Cells.Item(8, 17).Value = 7.99876828516
$ws.Cells.Item(8, 18).Value = 71.98891456644
$ws.Cells.Item(8, 19).Value = 0.1550577158270771
$ws.Cells.Item(8, 20).Value = 0.155057715827077

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.7674366666666667
$ws.Cells.Item(9, 8).Value = 2.30231
$ws.Cells.Item(9, 9).Value = 0.1816273560550486
$ws.Cells.Item(9, 10).Value = 0.1816273560550485
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.048156333333333
$ws.Cells.Item(9, 14).Value = 3.144469
$ws.Cells.Item(9, 15).Value = 0.08585342255879412
$ws.Cells.Item(9, 16).Value = 0.08585342255879413
$ws.Cells.Item(9, 17).Value = 0.8043936025988888
$ws.Cells.Item(9, 18).Value = 7.23954242339
$ws.Cells.Item(9, 19).Value = 0.01559333014763064
$ws.Cells.Item(9, 20).Value = 0.01559333014763064

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.7674366666666667
$ws.Cells.Item(10, 8).Value = 2.30231
$ws.Cells.Item(10, 9).Value = 0.1816273560550486
$ws.Cells.Item(10, 10).Value = 0.1816273560550485
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.7378083333333333
$ws.Cells.Item(10, 14).Value = 2.213425
$ws.Cells.Item(10, 15).Value = 0.06043313253436396
$ws.Cells.Item(10, 16).Value = 0.06043313253436396
$ws.Cells.Item(10, 17).Value = 0.5662211679722222
$ws.Cells.Item(10, 18).Value = 5.095990511749999
$ws.Cells.Item(10, 19).Value = 0.01097631008034086
$ws.Cells.Item(10, 20).Value = 0.01097631008034086
